$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 81. This shifts the existing rows 81-86
# (the historical "Coco" price records) down to rows 82-87, matching the
# weekly update pattern: a new record is added at the top of the block
# and the previous ones slide down.
$ws.Rows.Item(81).Insert()

# Populate the new row 81 with this week's record (copy of the row's
# static/category columns, with the new date / volume / price figures).
$ws.Cells.Item(81, 1).Value = 10
$ws.Cells.Item(81, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(81, 3).Value = "La Araucanía"
$ws.Cells.Item(81, 4).Value = 44769
$ws.Cells.Item(81, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(81, 5).Value = 9
$ws.Cells.Item(81, 6).Value = "Fruta"
$ws.Cells.Item(81, 7).Value = 100108
$ws.Cells.Item(81, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(81, 9).Value = 100108007
$ws.Cells.Item(81, 10).Value = "Coco"
$ws.Cells.Item(81, 11).Value = "Sin especificar"
$ws.Cells.Item(81, 12).Value = "Primera"
$ws.Cells.Item(81, 13).Value = 20
$ws.Cells.Item(81, 14).Value = 30000
$ws.Cells.Item(81, 15).Value = 30000
$ws.Cells.Item(81, 16).Value = 30000
$ws.Cells.Item(81, 17).Value = "`$/malla 20 unidades"
$ws.Cells.Item(81, 18).Value = "Perú"
$ws.Cells.Item(81, 19).Value = 1500
$ws.Cells.Item(81, 20).Value = 20
